# Update cryptos list (Tue Sep  5 15:14:35 UTC 2023) with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.863.94"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.42"
$ws.Range("E3").Value = "  +0.73%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.68"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5062"
$ws.Range("E6").Value = "  +0.47%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2586"
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06439"
$ws.Range("E9").Value = "  +1.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  +4.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07814"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.300"
$ws.Range("E12").Value = "  +1.55%  "

# Rows 13-14: WrappedEther and WrappedliquidstakedEther2.0 swap places
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.15"
$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.634.05"
$ws.Range("E14").Value = "  +0.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5657"
$ws.Range("E15").Value = "  +3.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7650"
$ws.Range("E16").Value = "  -0.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("E17").Value = "  -0.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.873.74"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("E19").Value = "  +0.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.08"
$ws.Range("E20").Value = "  +0.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.368"
$ws.Range("E21").Value = "  -1.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.952"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.129"
$ws.Range("E23").Value = "  +1.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.794"
$ws.Range("E25").Value = "  -6.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.31"
$ws.Range("E26").Value = "  -1.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1254"
$ws.Range("E27").Value = "  +1.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.852"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.244"
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04950"
$ws.Range("E31").Value = "  +1.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.319"
$ws.Range("E32").Value = "  +2.60%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.255"
$ws.Range("E33").Value = "  +2.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.581"
$ws.Range("E34").Value = "  +2.90%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.378"
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9089"
$ws.Range("E36").Value = "  +1.58%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.573"
$ws.Range("E37").Value = "  +1.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5555"
$ws.Range("E38").Value = "  +1.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.127.18"
$ws.Range("E39").Value = "  +0.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01563"
$ws.Range("E40").Value = "  +0.87%  "

# Row 41
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.529"
$ws.Range("E42").Value = "  -0.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8011"
$ws.Range("E43").Value = "  -0.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.25"
$ws.Range("E44").Value = "  +1.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.774.00"
$ws.Range("E45").Value = "  +0.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -5.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.71"
$ws.Range("E47").Value = "  +1.93%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4279"
$ws.Range("E48").Value = "  -3.96%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.783"
$ws.Range("E49").Value = "  +3.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05023"
$ws.Range("E50").Value = "  -2.33%  "

# Row 51
$ws.Range("E51").Value = "  +0.13%  "

